$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5277743333333333
$ws.Range("H2").Value = 1.583323
$ws.Range("I2").Value = 0.01400965908295571
$ws.Range("J2").Value = 0.01400965908295571
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.03057133333333333
$ws.Range("N2").Value = 0.09171399999999999
$ws.Range("O2").Value = 0.1084248755136686
$ws.Range("P2").Value = 0.1084248755136687
$ws.Range("Q2").Value = 0.01613476506911111
$ws.Range("R2").Value = 0.145212885622
$ws.Range("S2").Value = 0.00151899554205841
$ws.Range("T2").Value = 0.00151899554205841

$ws.Range("G3").Value = 0.5277743333333333
$ws.Range("H3").Value = 1.583323
$ws.Range("I3").Value = 0.01400965908295571
$ws.Range("J3").Value = 0.01400965908295571
$ws.Range("O3").Value = 0.8915751244863314
$ws.Range("P3").Value = 0.8915751244863314
$ws.Range("Q3").Value = 0.1326757822584445
$ws.Range("R3").Value = 1.194082040326
$ws.Range("S3").Value = 0.0124906635408973
$ws.Range("T3").Value = 0.0124906635408973

$ws.Range("I4").Value = 0.9003255417707673
$ws.Range("J4").Value = 0.9003255417707672
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.03057133333333333
$ws.Range("N4").Value = 0.09171399999999999
$ws.Range("O4").Value = 0.1084248755136686
$ws.Range("P4").Value = 0.1084248755136687
$ws.Range("Q4").Value = 1.036894689312222
$ws.Range("R4").Value = 9.332052203809999
$ws.Range("S4").Value = 0.09761768478827172
$ws.Range("T4").Value = 0.09761768478827174

$ws.Range("I5").Value = 0.9003255417707673
$ws.Range("J5").Value = 0.9003255417707672
$ws.Range("O5").Value = 0.8915751244863314
$ws.Range("P5").Value = 0.8915751244863314
$ws.Range("S5").Value = 0.8027078569824956
$ws.Range("T5").Value = 0.8027078569824955

$ws.Range("I6").Value = 0.08566479914627706
$ws.Range("J6").Value = 0.08566479914627706
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03057133333333333
$ws.Range("N6").Value = 0.09171399999999999
$ws.Range("O6").Value = 0.1084248755136686
$ws.Range("P6").Value = 0.1084248755136687
$ws.Range("Q6").Value = 0.09865917512577776
$ws.Range("R6").Value = 0.8879325761319998
$ws.Range("S6").Value = 0.009288195183338517
$ws.Range("T6").Value = 0.009288195183338521

$ws.Range("I7").Value = 0.08566479914627706
$ws.Range("J7").Value = 0.08566479914627706
$ws.Range("O7").Value = 0.8915751244863314
$ws.Range("P7").Value = 0.8915751244863314
$ws.Range("S7").Value = 0.07637660396293854
$ws.Range("T7").Value = 0.07637660396293854
